# hours update source add
# Fill in the previously-blank rows 45 and 46 on Sheet1 with new log entries
# and move the active selection to A47 (the next empty row).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 45: 3/23/2010, 2.5 hours, "Group Meeting"
$ws.Range("A45").Value = [DateTime]"2010-03-23"
$ws.Range("B45").Value = 2.5
$ws.Range("C45").Value = "Group Meeting"

# Row 46: 3/23/2010, 1 hour, "Weekly Meeting"
$ws.Range("A46").Value = [DateTime]"2010-03-23"
$ws.Range("B46").Value = 1
$ws.Range("C46").Value = "Weekly Meeting"

# Move selection to the next empty row, as in the saved workbook state.
$ws.Range("A47").Select()
